$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# --- Update shared-string (text) columns ---
# Order matters: the engine rebuilds the shared-string table in the order
# cells are assigned, so we assign column-by-column (B, then D, then E)
# to reproduce the original table's column-major ordering.

# best_params (column B)
$ws.Range("B2").Value = "{'max_depth': 5, 'min_samples_leaf': 2, 'min_samples_split': 5}"
$ws.Range("B3").Value = "{'max_depth': 20, 'min_samples_split': 2, 'n_estimators': 300}"
$ws.Range("B4").Value = "{'learning_rate': 0.05, 'n_estimators': 150, 'num_leaves': 31}"
$ws.Range("B5").Value = "{'learning_rate': 0.1, 'max_depth': 5, 'n_estimators': 250}"

# best_model (column D)
$ws.Range("D2").Value = "DecisionTreeClassifier(max_depth=5, min_samples_leaf=2, min_samples_split=5)"
$ws.Range("D3").Value = "RandomForestClassifier(max_depth=20, n_estimators=300)"
$ws.Range("D4").Value = "LGBMClassifier(learning_rate=0.05, n_estimators=150)"
$ws.Range("D5").Value = "XGBClassifier(base_score=None, booster=None, callbacks=None," + $nl + "              colsample_bylevel=None, colsample_bynode=None," + $nl + "              colsample_bytree=None, device=None, early_stopping_rounds=None," + $nl + "              enable_categorical=True, eval_metric=None, feature_types=None," + $nl + "              gamma=None, grow_policy=None, importance_type=None," + $nl + "              interaction_constraints=None, learning_rate=0.1, max_bin=None," + $nl + "              max_cat_threshold=None, max_cat_to_onehot=None," + $nl + "              max_delta_step=None, max_depth=5, max_leaves=None," + $nl + "              min_child_weight=None, missing=nan, monotone_constraints=None," + $nl + "              multi_strategy=None, n_estimators=250, n_jobs=None," + $nl + "              num_parallel_tree=None, random_state=None, ...)"

# confusion_matrix (column E)
$ws.Range("E2").Value = "[[220  57]" + $nl + " [ 33 416]]"
$ws.Range("E3").Value = "[[219  58]" + $nl + " [ 26 423]]"
$ws.Range("E4").Value = "[[236  41]" + $nl + " [ 23 426]]"
$ws.Range("E5").Value = "[[232  45]" + $nl + " [ 30 419]]"

# --- Update numeric columns ---
# Row 2 (CART)
$ws.Range("C2").Value = 0.8701611565351346
$ws.Range("G2").Value = 57
$ws.Range("I2").Value = 220
$ws.Range("J2").Value = 0.8757048798997842
$ws.Range("K2").Value = 0.8760330578512396
$ws.Range("L2").Value = 0.8748397117703904
$ws.Range("M2").Value = 21.08915519714355

# Row 3 (Random Forest)
$ws.Range("C3").Value = 0.897023343998104
$ws.Range("F3").Value = 423
$ws.Range("G3").Value = 58
$ws.Range("H3").Value = 26
$ws.Range("I3").Value = 219
$ws.Range("J3").Value = 0.8849348615582383
$ws.Range("K3").Value = 0.8842975206611571
$ws.Range("L3").Value = 0.8827416648025406
$ws.Range("M3").Value = 676.3029205799103

# Row 4 (LightGBM)
$ws.Range("C4").Value = 0.8987617016234151
$ws.Range("F4").Value = 426
$ws.Range("G4").Value = 41
$ws.Range("H4").Value = 23
$ws.Range("I4").Value = 236
$ws.Range("J4").Value = 0.9118207220212281
$ws.Range("K4").Value = 0.9118457300275482
$ws.Range("L4").Value = 0.9112316723071997
$ws.Range("M4").Value = 197.9825568199158

# Row 5 (XGBoost)
$ws.Range("C5").Value = 0.8970375636923805
$ws.Range("F5").Value = 419
$ws.Range("G5").Value = 45
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = 232
$ws.Range("J5").Value = 0.896332209541876
$ws.Range("K5").Value = 0.8966942148760331
$ws.Range("L5").Value = 0.8961053739790885
$ws.Range("M5").Value = 263.7599172592163

# The multi-line text we just wrote triggers Excel's implicit row auto-height
# recalculation; restore rows to their original (no explicit height) state.
$ws.Rows("2:5").AutoFit()
